# Edit script for "Översikt BORLÄNGE.xlsx" — automatic update of files.
#
# Summary of changes (per the supplied diff):
#   1. Every row's "Förändrad" (C column) timestamp moves from 45189 -> 45190
#      (rows 2..301), and the same value is used for the newly appended row 302.
#   2. Rows 3 and 4 swap which case they describe:
#        - Row 3 used to hold "A 55504-2022"; it now holds "A 73613-2021"
#          (with a couple of updated counts and an extra species in the list).
#        - Row 4 used to hold "A 73613-2021"; it now holds "A 55504-2022"
#          (the data that used to live in row 3, content unchanged otherwise).
#   3. A brand-new case "A 44240-2023" is appended as row 302.
#   4. Row 301 keeps its data, just gets its row height re-stamped (cosmetic).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bulk date refresh: C2:C301 45189 -> 45190.
# ---------------------------------------------------------------------------
$ws.Range("C2:C301").Value = 45190

# ---------------------------------------------------------------------------
# 2) Swap the content of row 3 ("A 55504-2022") and row 4 ("A 73613-2021"),
#    writing the post-swap values (which differ slightly from a pure swap).
# ---------------------------------------------------------------------------

# --- New row 3: "A 73613-2021" --------------------------------------------
$ws.Range("A3").Value = "A 73613-2021"
$ws.Range("B3").Value = 44552
$ws.Range("C3").Value = 45190
$ws.Range("D3").Value = "DALARNAS LÄN"
$ws.Range("E3").Value = "BORLÄNGE"
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = 28
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 9
$ws.Range("R3").Value = "Garnlav`nMotaggsvamp`nSvart taggsvamp`nTallticka`nTretåig hackspett`nUllticka`nDropptaggsvamp`nRostfläck`nVedticka"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 73613-2021.xlsx", "A 73613-2021")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 73613-2021.png", "A 73613-2021")'
$ws.Range("U3").ClearContents()
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 73613-2021.docx", "A 73613-2021")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 73613-2021.docx", "A 73613-2021")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 73613-2021.docx", "A 73613-2021")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 73613-2021.docx", "A 73613-2021")'

# --- New row 4: "A 55504-2022" --------------------------------------------
$ws.Range("A4").Value = "A 55504-2022"
$ws.Range("B4").Value = 44887
$ws.Range("C4").Value = 45190
$ws.Range("D4").Value = "DALARNAS LÄN"
$ws.Range("E4").Value = "BORLÄNGE"
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = 15.1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 9
$ws.Range("R4").Value = "Knärot`nRynkskinn`nGammelgranskål`nGarnlav`nGranticka`nRosenticka`nUllticka`nSkuggblåslav`nVedticka"
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 55504-2022.xlsx", "A 55504-2022")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 55504-2022.png", "A 55504-2022")'
$ws.Range("U4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 55504-2022.png", "A 55504-2022")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 55504-2022.docx", "A 55504-2022")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 55504-2022.docx", "A 55504-2022")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 55504-2022.docx", "A 55504-2022")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 55504-2022.docx", "A 55504-2022")'

# Every row in this sheet is pinned to a 15pt row height (not auto-fit), even
# though R has wrapped, multi-line text. Re-write the content above can make
# the host auto-grow rows 3/4 to fit the (now 9-line) wrapped text, so pin
# the height back down to match the workbook's existing convention.
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15

# ---------------------------------------------------------------------------
# 3) Row 301 is unchanged in content; only re-stamp its row height so it
#    serialises with an explicit ht/customHeight (cosmetic, matches diff).
# ---------------------------------------------------------------------------
$ws.Rows.Item(301).RowHeight = 15

# ---------------------------------------------------------------------------
# 4) Append new row 302: "A 44240-2023".
# ---------------------------------------------------------------------------
$ws.Range("A302").Value = "A 44240-2023"
$ws.Range("B302").Value = 45188
$ws.Range("B302").NumberFormat = "YYYY-MM-DD"
$ws.Range("C302").Value = 45190
$ws.Range("C302").NumberFormat = "YYYY-MM-DD"
$ws.Range("D302").Value = "DALARNAS LÄN"
$ws.Range("E302").Value = "BORLÄNGE"
$ws.Range("G302").Value = 0.9
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = 0
$ws.Range("J302").Value = 0
$ws.Range("K302").Value = 0
$ws.Range("L302").Value = 0
$ws.Range("M302").Value = 0
$ws.Range("N302").Value = 0
$ws.Range("O302").Value = 0
$ws.Range("P302").Value = 0
$ws.Range("Q302").Value = 0
$ws.Range("R302").Value = ""
$ws.Range("R302").WrapText = $true
